$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value while preserving the cell as TEXT
# (the source workbook stores every data cell as an inline string, so a
# plain numeric assignment would be re-typed by Excel as a number).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    # Drop back to the default "Normal" style so we don't leave a stray
    # quote-prefix style behind (keeps formatting identical to before).
    $range.Style = "Normal"
}

# --- Simple price refreshes (column D) ---
Set-TextValue $ws.Range("D2")  "245.14"
Set-TextValue $ws.Range("D4")  "5.402"
Set-TextValue $ws.Range("D5")  "0.06046"
Set-TextValue $ws.Range("D6")  "3.395"
Set-TextValue $ws.Range("D7")  "0.8080"
Set-TextValue $ws.Range("D8")  "0.9328"
Set-TextValue $ws.Range("D9")  "0.1424"
Set-TextValue $ws.Range("D10") "0.07432"
Set-TextValue $ws.Range("D12") "0.03070"
Set-TextValue $ws.Range("D13") "0.09368"
Set-TextValue $ws.Range("D14") "3.935"
Set-TextValue $ws.Range("D15") "0.001597"
Set-TextValue $ws.Range("D16") "0.04837"
Set-TextValue $ws.Range("D17") "0.0005945"
Set-TextValue $ws.Range("D18") "0.005389"
Set-TextValue $ws.Range("D19") "0.004162"
Set-TextValue $ws.Range("D20") "0.0009862"
Set-TextValue $ws.Range("D22") "3.649"
Set-TextValue $ws.Range("D23") "6.443"
Set-TextValue $ws.Range("D26") "0.1296"
Set-TextValue $ws.Range("D27") "0.0002448"
Set-TextValue $ws.Range("D40") "0.03976"

# --- Rows 41-43: coin list shuffled by one position (KickToken moved to
#     the top, BKEXToken and CEJI shifted down), plus new prices ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.006384"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1074"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.002902"
$ws.Range("E43").Value = "42CEJICEJI"

# --- More simple price refreshes (column D) ---
Set-TextValue $ws.Range("D44") "0.005980"
Set-TextValue $ws.Range("D45") "0.00005167"
Set-TextValue $ws.Range("D46") "0.00000000751"
Set-TextValue $ws.Range("D47") "0.0005805"
Set-TextValue $ws.Range("D48") "0.8206"
